$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 106
$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 5000
$ws.Range("K106").Value = 5000
$ws.Range("M106").Value = -4369

# ALC row 111
$ws.Range("H111").Value = 6014
$ws.Range("I111").Value = 6014
$ws.Range("K111").Value = 18042
$ws.Range("M111").Value = -14975

# ALC row 112
$ws.Range("H112").Value = 113676.164
$ws.Range("J112").Value = 74298
$ws.Range("L112").Value = 222894
$ws.Range("N112").Value = -225110

# ALC row 125
$ws.Range("H125").Value = 28127222
$ws.Range("J125").Value = 37501964
$ws.Range("L125").Value = 337517676
$ws.Range("N125").Value = -337522596

# ALC row 127
$ws.Range("H127").Value = 5802
$ws.Range("I127").Value = 3550
$ws.Range("J127").Value = 7153.2
$ws.Range("K127").Value = 10650
$ws.Range("L127").Value = 21459.6
$ws.Range("M127").Value = -5690
$ws.Range("N127").Value = -31379.6

# ALC row 129
$ws.Range("H129").Value = 13831.714
$ws.Range("I129").Value = 1804.1666
$ws.Range("K129").Value = 5412.4998
$ws.Range("M129").Value = -412.4997999999996

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 1944.6875
$ws.Range("I32").Value = 1800.7301
$ws.Range("K32").Value = 1800.7301
$ws.Range("M32").Value = -1513.7301

# ARM row 61
$ws.Range("H61").Value = 125002690
$ws.Range("I61").Value = 333335170
$ws.Range("J61").Value = 3199.6
$ws.Range("K61").Value = 333335170
$ws.Range("L61").Value = 3199.6
$ws.Range("M61").Value = -333334958
$ws.Range("N61").Value = -3623.6

# ARM row 74
$ws.Range("H74").Value = 25645474
$ws.Range("I74").Value = 32261808
$ws.Range("J74").Value = 7187.25
$ws.Range("K74").Value = 32261808
$ws.Range("L74").Value = 7187.25
$ws.Range("M74").Value = -32260934
$ws.Range("N74").Value = -8935.25

# ARM row 77
$ws.Range("H77").Value = 25645474
$ws.Range("I77").Value = 32261808
$ws.Range("J77").Value = 7187.25
$ws.Range("K77").Value = 161309040
$ws.Range("L77").Value = 35936.25
$ws.Range("M77").Value = -161304672
$ws.Range("N77").Value = -44672.25

# ARM row 80
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5000
$ws.Range("N80").ClearContents()
$ws.Range("M80").Value = -4002

# ARM row 83
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("N83").ClearContents()
$ws.Range("M83").Value = -10008

# ARM row 97
$ws.Range("H97").Value = 385.125
$ws.Range("I97").Value = 385.125
$ws.Range("K97").Value = 385.125
$ws.Range("M97").Value = 110.875

# ARM row 101
$ws.Range("H101").Value = 162257.14
$ws.Range("J101").Value = 162257.14
$ws.Range("L101").Value = 162257.14
$ws.Range("N101").Value = -168747.14

# ARM row 110
$ws.Range("H110").Value = 46894.184
$ws.Range("I110").Value = 51473.25
$ws.Range("K110").Value = 51473.25
$ws.Range("M110").Value = -49428.25

# ARM row 132
$ws.Range("H132").Value = 3335511.8
$ws.Range("I132").Value = 1925303.2
$ws.Range("J132").Value = 18001680
$ws.Range("K132").Value = 5775909.6
$ws.Range("L132").Value = 54005040
$ws.Range("M132").Value = -5773379.6
$ws.Range("N132").Value = -54010100

# ARM row 136
$ws.Range("H136").Value = 125002690
$ws.Range("I136").Value = 333335170
$ws.Range("J136").Value = 3199.6
$ws.Range("K136").Value = 1000005510
$ws.Range("L136").Value = 9598.799999999999
$ws.Range("M136").Value = -1000002960
$ws.Range("N136").Value = -14698.8

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 2018.8334
$ws.Range("I20").Value = 1915.875
$ws.Range("K20").Value = 1915.875
$ws.Range("M20").Value = -1668.875

# BSM row 107
$ws.Range("H107").Value = 41026.47
$ws.Range("I107").Value = 1766.3334
$ws.Range("J107").Value = 135250.8
$ws.Range("K107").Value = 1766.3334
$ws.Range("L107").Value = 135250.8
$ws.Range("M107").Value = 153.6666
$ws.Range("N107").Value = -139090.8

$ws = $wb.Worksheets.Item("CRP")
# CRP row 15
$ws.Range("H15").Value = 1631.25
$ws.Range("J15").Value = 1631.25
$ws.Range("L15").Value = 1631.25
$ws.Range("N15").Value = -1971.25

$ws = $wb.Worksheets.Item("CUL")
# CUL row 12
$ws.Range("H12").Value = 363.33334
$ws.Range("I12").Value = 66
$ws.Range("J12").Value = 448.2857
$ws.Range("K12").Value = 198
$ws.Range("L12").Value = 1344.8571
$ws.Range("M12").Value = -25
$ws.Range("N12").Value = -1690.8571

# CUL row 43
$ws.Range("H43").Value = 32000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 32000
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -96228

# CUL row 107
$ws.Range("H107").Value = 3100
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 18000
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -21840

# CUL row 109
$ws.Range("H109").Value = 2063.3333
$ws.Range("I109").Value = 1595
$ws.Range("K109").Value = 4785
$ws.Range("M109").Value = -3745

# CUL row 131
$ws.Range("H131").Value = 1635
$ws.Range("I131").Value = 1268.75
$ws.Range("K131").Value = 3806.25
$ws.Range("M131").Value = 1233.75

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 3951.4443
$ws.Range("I80").Value = 3801.3076
$ws.Range("K80").Value = 3801.3076
$ws.Range("M80").Value = -2803.3076

# GSM row 83
$ws.Range("H83").Value = 3951.4443
$ws.Range("I83").Value = 3801.3076
$ws.Range("K83").Value = 19006.538
$ws.Range("M83").Value = -14014.538

# GSM row 97
$ws.Range("H97").Value = 1700
$ws.Range("I97").Value = 200
$ws.Range("K97").Value = 200
$ws.Range("M97").Value = 296

# GSM row 102
$ws.Range("H102").Value = 6370.1763
$ws.Range("I102").Value = 3716.2
$ws.Range("K102").Value = 3716.2
$ws.Range("M102").Value = -2094.2

# GSM row 111
$ws.Range("H111").Value = 48333.332

# GSM row 113
$ws.Range("H113").Value = 33784.324
$ws.Range("I113").Value = 42841.965
$ws.Range("J113").Value = 5605
$ws.Range("K113").Value = 42841.965
$ws.Range("L113").Value = 5605
$ws.Range("M113").Value = -40671.965
$ws.Range("N113").Value = -9945

# GSM row 122
$ws.Range("H122").Value = 7581.273
$ws.Range("I122").Value = 2949.6667
$ws.Range("K122").Value = 8849.000100000001
$ws.Range("M122").Value = -6399.000100000001

# GSM row 123
$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -104900

# GSM row 132
$ws.Range("H132").Value = 4645602.5
$ws.Range("I132").Value = 5015031
$ws.Range("K132").Value = 15045093
$ws.Range("M132").Value = -15042563

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 2094.2222
$ws.Range("I7").Value = 2107.25
$ws.Range("K7").Value = 2107.25
$ws.Range("M7").Value = -1995.25

# LTW row 22
$ws.Range("H22").Value = 2659.25
$ws.Range("I22").Value = 2742.6365
$ws.Range("J22").Value = 2475.8
$ws.Range("K22").Value = 2742.6365
$ws.Range("L22").Value = 2475.8
$ws.Range("M22").Value = -2447.6365
$ws.Range("N22").Value = -3065.8

# LTW row 27
$ws.Range("H27").Value = 2659.25
$ws.Range("I27").Value = 2742.6365
$ws.Range("J27").Value = 2475.8
$ws.Range("K27").Value = 2742.6365
$ws.Range("L27").Value = 2475.8
$ws.Range("M27").Value = -2635.6365
$ws.Range("N27").Value = -2689.8

# LTW row 46
$ws.Range("H46").Value = 2716.6667
$ws.Range("I46").Value = 2390
$ws.Range("K46").Value = 2390
$ws.Range("M46").Value = -2202

# LTW row 55
$ws.Range("H55").Value = 685.2222
$ws.Range("I55").Value = 288.4
$ws.Range("K55").Value = 288.4
$ws.Range("M55").Value = -115.4

# LTW row 122
$ws.Range("H122").Value = 3190.75
$ws.Range("I122").Value = 3172.25
$ws.Range("J122").Value = 3246.25
$ws.Range("K122").Value = 9516.75
$ws.Range("L122").Value = 9738.75
$ws.Range("M122").Value = -7066.75
$ws.Range("N122").Value = -14638.75

# LTW row 126
$ws.Range("H126").Value = 2094.2222
$ws.Range("I126").Value = 2107.25
$ws.Range("K126").Value = 6321.75
$ws.Range("M126").Value = -3851.75

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 1629.9445
$ws.Range("I122").Value = 1611.6666
$ws.Range("J122").Value = 1721.3334
$ws.Range("K122").Value = 4834.9998
$ws.Range("L122").Value = 5164.0002
$ws.Range("M122").Value = -2384.9998
$ws.Range("N122").Value = -10475.9998

Write-Output "Applied all Spriggan_Profits updates"